# Update the Confidence parameter (delta Pn) on the "Data" sheet.
# Cell H3 drives the confidence-interval calculations (columns B:E) that are
# plotted in the two charts on the "Chart" sheet. Changing H3 from 0.2 to
# 0.74 propagates automatically through all the dependent formulas and the
# charts' cached values when Excel recalculates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

$ws.Range("H3").Value = 0.74

# Reflect the cell that was last selected/active on the "Data" sheet.
$ws.Range("C13").Select()

$excel.CalculateFullRebuild()
